# Apply weekly update of fruit/vegetable prices (Haba) by reshuffling
# the D, J, K, L, M, O, P columns across rows 2-13 as described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row: Fecha(D), Volumen(J), Precio minimo(K), Precio maximo(L),
# Precio promedio ponderado(M), Origen(O), Precio $/Kg(P)
$rows = @{
  2  = @{ D = 44453; J = 55; K = 14000; L = 15000; M = 14455; O = "Provincia de Limarí"; P = 578 }
  3  = @{ D = 44425; J = 25; K = 14000; L = 14000; M = 14000; O = "Provincia de Limarí"; P = 560 }
  4  = @{ D = 44418; J = 12; K = 15000; L = 15000; M = 15000; O = "Provincia de Limarí"; P = 600 }
  5  = @{ D = 44432; J = 15; K = 14000; L = 14000; M = 14000; O = "Provincia del Elquí"; P = 560 }
  6  = @{ D = 44435; J = 15; K = 14000; L = 14000; M = 14000; O = "Provincia de Limarí"; P = 560 }
  7  = @{ D = 44435; J = 15; K = 14000; L = 14000; M = 14000; O = "Provincia del Elquí"; P = 560 }
  8  = @{ D = 44449; J = 30; K = 16000; L = 16000; M = 16000; O = "Provincia de Limarí"; P = 640 }
  9  = @{ D = 44421; J = 20; K = 15000; L = 15000; M = 15000; O = "Provincia de Limarí"; P = 600 }
  11 = @{ D = 44446; J = 15; K = 13000; L = 13000; M = 13000; O = "Provincia de Limarí"; P = 520 }
  12 = @{ D = 44467; J = 35; K = 12000; L = 12000; M = 12000; O = "Provincia de Limarí"; P = 480 }
  13 = @{ D = 44340; J = 25; K = 15000; L = 15000; M = 15000; O = "Provincia de Limarí"; P = 600 }
}

foreach ($r in $rows.Keys) {
  $v = $rows[$r]
  $ws.Cells.Item($r, 4).Value2  = $v.D   # D: Fecha
  $ws.Cells.Item($r, 10).Value2 = $v.J   # J: Volumen
  $ws.Cells.Item($r, 11).Value2 = $v.K   # K: Precio minimo
  $ws.Cells.Item($r, 12).Value2 = $v.L   # L: Precio maximo
  $ws.Cells.Item($r, 13).Value2 = $v.M   # M: Precio promedio ponderado
  $ws.Cells.Item($r, 15).Value  = $v.O   # O: Origen
  $ws.Cells.Item($r, 16).Value2 = $v.P   # P: Precio $/Kg
}
